$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
}

Set-TextValue 2 4 "330.57"
Set-TextValue 2 5 "0.76%"
Set-TextValue 3 4 "39.69"
Set-TextValue 3 5 "-1.21%"
Set-TextValue 4 4 "5.690"
Set-TextValue 4 5 "0.42%"
Set-TextValue 5 4 "0.08070"
Set-TextValue 5 5 "-0.68%"
Set-TextValue 6 4 "8.635"
Set-TextValue 6 5 "-0.69%"
Set-TextValue 7 4 "4.476"
Set-TextValue 7 5 "-1.75%"
Set-TextValue 8 4 "1.958"
Set-TextValue 8 5 "-0.72%"
Set-TextValue 9 4 "2.970"
Set-TextValue 9 5 "-0.24%"
Set-TextValue 10 4 "0.9234"
Set-TextValue 10 5 "-2.74%"
Set-TextValue 11 4 "0.1248"
Set-TextValue 11 5 "-1.34%"
Set-TextValue 12 4 "0.1950"
Set-TextValue 12 5 "-2.31%"
Set-TextValue 13 4 "8.698"
Set-TextValue 13 5 "15.88%"
Set-TextValue 14 4 "0.09237"
Set-TextValue 14 5 "0.39%"
Set-TextValue 15 4 "0.03597"
Set-TextValue 15 5 "1.29%"
Set-TextValue 16 4 "0.1051"
Set-TextValue 16 5 "9.32%"
Set-TextValue 17 4 "0.001307"
Set-TextValue 17 5 "-0.22%"
Set-TextValue 18 4 "0.006382"
Set-TextValue 18 5 "5.08%"
Set-TextValue 19 5 "-0.18%"
Set-TextValue 20 4 "0.3459"
Set-TextValue 20 5 "-1.60%"
Set-TextValue 21 4 "0.1370"
Set-TextValue 21 5 "-2.16%"
Set-TextValue 22 5 "4.23%"
Set-TextValue 23 4 "0.04429"
Set-TextValue 23 5 "-0.34%"
Set-TextValue 24 5 "-0.02%"
Set-TextValue 25 4 "0.004490"
Set-TextValue 25 5 "3.50%"
Set-TextValue 26 4 "0.0001199"
Set-TextValue 26 5 "0.66%"
Set-TextValue 39 4 "0.02582"
Set-TextValue 39 5 "2.15%"
Set-TextValue 40 4 "0.05484"
Set-TextValue 40 5 "5.44%"
Set-TextValue 41 4 "0.007532"
Set-TextValue 41 5 "-3.21%"
Set-TextValue 42 4 "0.009888"
Set-TextValue 42 5 "9.19%"
Set-TextValue 43 4 "0.1412"
Set-TextValue 43 5 "-1.14%"
Set-TextValue 44 4 "0.002105"
Set-TextValue 44 5 "-3.96%"
Set-TextValue 45 5 "8.82%"
Set-TextValue 46 4 "0.00006801"
Set-TextValue 46 5 "1.04%"
Set-TextValue 47 4 "0.00000000749"
Set-TextValue 47 5 "-0.24%"
Set-TextValue 48 4 "0.003049"
Set-TextValue 48 5 "6.08%"
Set-TextValue 49 4 "0.002279"
Set-TextValue 49 5 "-0.91%"
Set-TextValue 50 4 "0.00002097"
Set-TextValue 50 5 "-0.24%"
Set-TextValue 51 4 "0.0001997"
Set-TextValue 51 5 "-0.24%"
